# Scheduled-runner update: refresh ALC/ARM/CRP/CUL/GSM/LTW/WVR leve-profit sheets
# with latest Universalis market-price snapshots (currentAveragePrice* / Leve* columns).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 22224650
$ws.Range("I15").Value = 22224650
$ws.Range("K15").Value = 66673950
$ws.Range("M15").Value = -66673781
$ws.Range("H62").Value = 4280.7144
$ws.Range("I62").Value = 2950
$ws.Range("J62").Value = 4813
$ws.Range("K62").Value = 2950
$ws.Range("L62").Value = 4813
$ws.Range("M62").Value = -2326
$ws.Range("N62").Value = -6061
$ws.Range("H65").Value = 4280.7144
$ws.Range("I65").Value = 2950
$ws.Range("J65").Value = 4813
$ws.Range("K65").Value = 14750
$ws.Range("L65").Value = 24065
$ws.Range("M65").Value = -11630
$ws.Range("N65").Value = -30305
$ws.Range("H125").Value = 2079.2307
$ws.Range("I125").Value = 1872
$ws.Range("J125").Value = 2256.8572
$ws.Range("K125").Value = 16848
$ws.Range("L125").Value = 20311.7148
$ws.Range("M125").Value = -14388
$ws.Range("N125").Value = -25231.7148
$ws.Range("H132").Value = 4256819
$ws.Range("I132").Value = 4652155.5
$ws.Range("J132").Value = 6950
$ws.Range("K132").Value = 13956466.5
$ws.Range("L132").Value = 20850
$ws.Range("M132").Value = -13953936.5
$ws.Range("N132").Value = -25910
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -40120
$ws.Range("H138").Value = 3458.913
$ws.Range("J138").Value = 6284.467
$ws.Range("L138").Value = 18853.401
$ws.Range("N138").Value = -29133.401

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3293.371
$ws.Range("I32").Value = 2193.518
$ws.Range("J32").Value = 9813.929
$ws.Range("K32").Value = 2193.518
$ws.Range("L32").Value = 9813.929
$ws.Range("M32").Value = -1906.518
$ws.Range("N32").Value = -10387.929
$ws.Range("H45").Value = 1609.3871
$ws.Range("I45").Value = 865.087
$ws.Range("J45").Value = 3749.25
$ws.Range("K45").Value = 865.087
$ws.Range("L45").Value = 3749.25
$ws.Range("M45").Value = -488.087
$ws.Range("N45").Value = -4503.25
$ws.Range("H74").Value = 926.86365
$ws.Range("I74").Value = 824.5
$ws.Range("K74").Value = 824.5
$ws.Range("M74").Value = 49.5
$ws.Range("H77").Value = 926.86365
$ws.Range("I77").Value = 824.5
$ws.Range("K77").Value = 4122.5
$ws.Range("M77").Value = 245.5
$ws.Range("H122").Value = 3091.5334
$ws.Range("I122").Value = 1435.9
$ws.Range("K122").Value = 4307.700000000001
$ws.Range("M122").Value = -1857.700000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2728.4285
$ws.Range("I16").Value = 699.6667
$ws.Range("J16").Value = 4250
$ws.Range("K16").Value = 699.6667
$ws.Range("L16").Value = 4250
$ws.Range("M16").Value = -412.6667
$ws.Range("N16").Value = -4824
$ws.Range("H31").Value = 2002584.9
$ws.Range("J31").Value = 3656.261
$ws.Range("L31").Value = 3656.261
$ws.Range("N31").Value = -4246.261
$ws.Range("H34").Value = 2002584.9
$ws.Range("J34").Value = 3656.261
$ws.Range("L34").Value = 3656.261
$ws.Range("N34").Value = -4060.261
$ws.Range("H86").Value = 4165
$ws.Range("I86").Value = 3556.3333
$ws.Range("J86").Value = 4556.2856
$ws.Range("K86").Value = 3556.3333
$ws.Range("L86").Value = 4556.2856
$ws.Range("M86").Value = -2433.3333
$ws.Range("N86").Value = -6802.2856
$ws.Range("H89").Value = 4165
$ws.Range("I89").Value = 3556.3333
$ws.Range("J89").Value = 4556.2856
$ws.Range("K89").Value = 17781.6665
$ws.Range("L89").Value = 22781.428
$ws.Range("M89").Value = -12165.6665
$ws.Range("N89").Value = -34013.428
$ws.Range("H113").Value = 2728.4285
$ws.Range("I113").Value = 699.6667
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 699.6667
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = 1470.3333
$ws.Range("N113").Value = -8590
$ws.Range("H122").Value = 2758.9048
$ws.Range("I122").Value = 2058.2856
$ws.Range("J122").Value = 4160.143
$ws.Range("K122").Value = 6174.8568
$ws.Range("L122").Value = 12480.429
$ws.Range("M122").Value = -3724.8568
$ws.Range("N122").Value = -17380.429
$ws.Range("H132").Value = 3915.0715
$ws.Range("I132").Value = 2847.6365
$ws.Range("K132").Value = 8542.9095
$ws.Range("M132").Value = -6012.9095
$ws.Range("H134").Value = 1853.9062
$ws.Range("I134").Value = 1458.0435
$ws.Range("K134").Value = 4374.1305
$ws.Range("M134").Value = -1839.1305

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2408
$ws.Range("I68").Value = 739.86365
$ws.Range("J68").Value = 3937.125
$ws.Range("K68").Value = 2219.59095
$ws.Range("L68").Value = 11811.375
$ws.Range("M68").Value = -1408.59095
$ws.Range("N68").Value = -13433.375
$ws.Range("H71").Value = 2408
$ws.Range("I71").Value = 739.86365
$ws.Range("J71").Value = 3937.125
$ws.Range("K71").Value = 6658.77285
$ws.Range("L71").Value = 35434.125
$ws.Range("M71").Value = -2602.77285
$ws.Range("N71").Value = -43546.125
$ws.Range("H122").Value = 868.8946999999999
$ws.Range("J122").Value = 1104.9166
$ws.Range("L122").Value = 9944.249400000001
$ws.Range("N122").Value = -14844.2494

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3459.2222
$ws.Range("I80").Value = 3180.5833
$ws.Range("J80").Value = 4016.5
$ws.Range("K80").Value = 3180.5833
$ws.Range("L80").Value = 4016.5
$ws.Range("M80").Value = -2182.5833
$ws.Range("N80").Value = -6012.5
$ws.Range("H83").Value = 3459.2222
$ws.Range("I83").Value = 3180.5833
$ws.Range("J83").Value = 4016.5
$ws.Range("K83").Value = 15902.9165
$ws.Range("L83").Value = 20082.5
$ws.Range("M83").Value = -10910.9165
$ws.Range("N83").Value = -30066.5
$ws.Range("H97").Value = 1425
$ws.Range("J97").Value = 1613.5555
$ws.Range("L97").Value = 1613.5555
$ws.Range("N97").Value = -2605.5555
$ws.Range("H122").Value = 5549.722
$ws.Range("I122").Value = 5109.6665
$ws.Range("J122").Value = 5989.778
$ws.Range("K122").Value = 15328.9995
$ws.Range("L122").Value = 17969.334
$ws.Range("M122").Value = -12878.9995
$ws.Range("N122").Value = -22869.334
$ws.Range("H126").Value = 2586.697
$ws.Range("I126").Value = 1401.1333
$ws.Range("J126").Value = 3574.6667
$ws.Range("K126").Value = 4203.3999
$ws.Range("L126").Value = 10724.0001
$ws.Range("M126").Value = -1733.3999
$ws.Range("N126").Value = -15664.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2546.4546
$ws.Range("I7").Value = 1642.4
$ws.Range("J7").Value = 3299.8333
$ws.Range("K7").Value = 1642.4
$ws.Range("L7").Value = 3299.8333
$ws.Range("M7").Value = -1530.4
$ws.Range("N7").Value = -3523.8333
$ws.Range("H36").Value = 30000
$ws.Range("J36").Value = 30000
$ws.Range("L36").Value = 30000
$ws.Range("N36").Value = -31124
$ws.Range("H122").Value = 2973.6316
$ws.Range("I122").Value = 2568.8965
$ws.Range("K122").Value = 7706.689499999999
$ws.Range("M122").Value = -5256.689499999999
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("N123").Value = 0
$ws.Range("H126").Value = 2546.4546
$ws.Range("I126").Value = 1642.4
$ws.Range("J126").Value = 3299.8333
$ws.Range("K126").Value = 4927.200000000001
$ws.Range("L126").Value = 9899.499899999999
$ws.Range("M126").Value = -2457.200000000001
$ws.Range("N126").Value = -14839.4999
$ws.Range("H132").Value = 5777.5557
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 6666.3335
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 19999.0005
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -25059.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2528.0715
$ws.Range("I107").Value = 874.75
$ws.Range("J107").Value = 3189.4
$ws.Range("K107").Value = 2624.25
$ws.Range("L107").Value = 9568.200000000001
$ws.Range("M107").Value = -704.25
$ws.Range("N107").Value = -13408.2
$ws.Range("H122").Value = 419057.12
$ws.Range("I122").Value = 626605.1
$ws.Range("J122").Value = 3961.125
$ws.Range("K122").Value = 1879815.3
$ws.Range("L122").Value = 11883.375
$ws.Range("M122").Value = -1877365.3
$ws.Range("N122").Value = -16783.375
$ws.Range("H126").Value = 2382638.2
$ws.Range("I126").Value = 1289.3572
$ws.Range("J126").Value = 7145336
$ws.Range("K126").Value = 3868.0716
$ws.Range("L126").Value = 21436008
$ws.Range("M126").Value = -1398.0716
$ws.Range("N126").Value = -21440948
$ws.Range("H132").Value = 377143.6
$ws.Range("J132").Value = 12139.417
$ws.Range("L132").Value = 36418.251
$ws.Range("N132").Value = -41478.251

Write-Output "Applied 224 cell updates across 7 sheets."
